{"js": "// Remove the trailing \"Ver no Jupiter ...\" / footer copyright block that\n// followed the \"LOQ4205: ...\" requirements paragraph, along with the blank\n// paragraph separating them, while leaving the rest of the document\n// (including the final blank + page-break paragraphs) untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the three consecutive paragraphs to delete:\n//   1) the blank paragraph right after \"LOQ4205: ...\"\n//   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3) the \"\u00a9 2020 ...\" copyright paragraph\nlet reqIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOQ4205\") !== -1) {\n    reqIndex = i;\n    break;\n  }\n}\n\nif (reqIndex !== -1 &&\n    items[reqIndex + 1] && items[reqIndex + 1].text.trim() === \"\" &&\n    items[reqIndex + 2] && items[reqIndex + 2].text.indexOf(\"Ver no Jupiter\") !== -1 &&\n    items[reqIndex + 3] && items[reqIndex + 3].text.indexOf(\"\\u00A9 2020\") !== -1) {\n  // Delete from last to first so the earlier deletes don't shift later ones.\n  items[reqIndex + 3].delete();\n  items[reqIndex + 2].delete();\n  items[reqIndex + 1].delete();\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / footer copyright block that\n# followed the \"LOQ4205: ...\" requirements paragraph, along with the blank\n# paragraph separating them, while leaving the rest of the document\n# (including the final blank + page-break paragraphs) untouched.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$reqIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.Contains(\"LOQ4205\")) {\n        $reqIndex = $i\n        break\n    }\n}\n\nif ($reqIndex -gt 0) {\n    $i1 = $reqIndex + 1\n    $i2 = $reqIndex + 2\n    $i3 = $reqIndex + 3\n\n    $blank = $d.Paragraphs.Item($i1).Range.Text.Trim()\n    $jupiter = $d.Paragraphs.Item($i2).Range.Text\n    $copyright = $d.Paragraphs.Item($i3).Range.Text\n\n    if ($blank -eq \"\" -and $jupiter.Contains(\"Ver no Jupiter\") -and $copyright.Contains(\"2020\")) {\n        # Delete from the highest index down so earlier deletes don't shift\n        # the still-to-be-deleted paragraphs.\n        $d.Paragraphs.Item($i3).Range.Delete()\n        $d.Paragraphs.Item($i2).Range.Delete()\n        $d.Paragraphs.Item($i1).Range.Delete()\n    }\n}\n"}
